$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1250183.9
$ws.Range("I2").Value = 1428772.2
$ws.Range("K2").Value = 1428772.2
$ws.Range("M2").Value = -1428659.2

$ws.Range("H6").Value = 29412028
$ws.Range("I6").Value = 35714412
$ws.Range("K6").Value = 107143236
$ws.Range("M6").Value = -107143124

$ws.Range("H12").Value = 250.71428
$ws.Range("I12").Value = 276
$ws.Range("K12").Value = 276
$ws.Range("M12").Value = -106

$ws.Range("H17").Value = 2245.7334
$ws.Range("J17").Value = 2334.7856
$ws.Range("L17").Value = 7004.3568
$ws.Range("N17").Value = -7340.3568

$ws.Range("H19").Value = 13890514
$ws.Range("J19").Value = 1750.5
$ws.Range("L19").Value = 1750.5
$ws.Range("N19").Value = -2100.5

$ws.Range("H41").Value = 496.58334
$ws.Range("I41").Value = 599.8570999999999
$ws.Range("K41").Value = 599.8570999999999
$ws.Range("M41").Value = -159.8570999999999

$ws.Range("H53").Value = 208.91667
$ws.Range("I53").Value = 104
$ws.Range("J53").Value = 283.85715
$ws.Range("K53").Value = 104
$ws.Range("L53").Value = 283.85715
$ws.Range("M53").Value = 533
$ws.Range("N53").Value = -1557.85715

$ws.Range("H111").Value = 3547
$ws.Range("I111").Value = 3609.7144
$ws.Range("J111").Value = 3327.5
$ws.Range("K111").Value = 10829.1432
$ws.Range("L111").Value = 9982.5
$ws.Range("M111").Value = -7762.143199999999
$ws.Range("N111").Value = -16116.5

$ws.Range("H129").Value = 1666.1666
$ws.Range("I129").Value = 999.75
$ws.Range("J129").Value = 2999
$ws.Range("K129").Value = 2999.25
$ws.Range("L129").Value = 8997
$ws.Range("M129").Value = 2000.75
$ws.Range("N129").Value = -18997

$ws.Range("H138").Value = 5058358.5
$ws.Range("J138").Value = 11116257
$ws.Range("L138").Value = 33348771
$ws.Range("N138").Value = -33359051

$ws.Range("H141").Value = 10199
$ws.Range("I141").Value = 9898.9
$ws.Range("K141").Value = 29696.7
$ws.Range("M141").Value = -24516.7

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 616.55554
$ws.Range("I97").Value = 616.55554
$ws.Range("K97").Value = 616.55554
$ws.Range("M97").Value = -120.55554

$ws.Range("H110").Value = 1856.742
$ws.Range("I110").Value = 2133.68
$ws.Range("J110").Value = 702.8333
$ws.Range("K110").Value = 2133.68
$ws.Range("L110").Value = 702.8333
$ws.Range("M110").Value = -88.67999999999984
$ws.Range("N110").Value = -4792.8333

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 951.5
$ws.Range("J5").Value = 268.66666
$ws.Range("L5").Value = 268.66666
$ws.Range("N5").Value = -494.66666

$ws.Range("H11").Value = 461.125
$ws.Range("I11").Value = 647.5
$ws.Range("J11").Value = 274.75
$ws.Range("K11").Value = 647.5
$ws.Range("L11").Value = 274.75
$ws.Range("M11").Value = -507.5
$ws.Range("N11").Value = -554.75

$ws.Range("H20").Value = 6455.909
$ws.Range("I20").Value = 7877.1875
$ws.Range("J20").Value = 2665.8333
$ws.Range("K20").Value = 7877.1875
$ws.Range("L20").Value = 2665.8333
$ws.Range("M20").Value = -7630.1875
$ws.Range("N20").Value = -3159.8333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 3580.0908
$ws.Range("I105").Value = 3929.3333
$ws.Range("J105").Value = 3161
$ws.Range("K105").Value = 3929.3333
$ws.Range("L105").Value = 3161
$ws.Range("M105").Value = -2182.3333
$ws.Range("N105").Value = -6655

$ws.Range("H132").Value = 1970.8823
$ws.Range("I132").Value = 1931.5385
$ws.Range("K132").Value = 5794.6155
$ws.Range("M132").Value = -3264.6155

$ws.Range("H134").Value = 2432.5386
$ws.Range("I134").Value = 2511.3635
$ws.Range("K134").Value = 7534.0905
$ws.Range("M134").Value = -4999.0905

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 481.2353
$ws.Range("I2").Value = 162.28572
$ws.Range("J2").Value = 704.5
$ws.Range("K2").Value = 973.71432
$ws.Range("L2").Value = 4227
$ws.Range("M2").Value = -860.71432
$ws.Range("N2").Value = -4453

$ws.Range("H4").Value = 88.89
$ws.Range("I4").Value = 89.28283
$ws.Range("J4").Value = 50
$ws.Range("K4").Value = 267.84849
$ws.Range("L4").Value = 150
$ws.Range("M4").Value = -155.84849
$ws.Range("N4").Value = -374

$ws.Range("H10").Value = 167.83333
$ws.Range("I10").Value = 192.4
$ws.Range("K10").Value = 577.2
$ws.Range("M10").Value = -438.2

$ws.Range("H12").Value = 167.18182
$ws.Range("J12").Value = 186.28572
$ws.Range("L12").Value = 558.85716
$ws.Range("N12").Value = -904.85716

$ws.Range("H107").Value = 2653.125
$ws.Range("J107").Value = 1857.4
$ws.Range("L107").Value = 5572.200000000001
$ws.Range("N107").Value = -9412.200000000001

$ws.Range("H134").Value = 4653.0513
$ws.Range("I134").Value = 1498.1765
$ws.Range("K134").Value = 4494.529500000001
$ws.Range("M134").Value = 575.4704999999994

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 2054.7273
$ws.Range("I97").Value = 2183.8
$ws.Range("K97").Value = 2183.8
$ws.Range("M97").Value = -1687.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 2111.7693
$ws.Range("I55").Value = 1307.25
$ws.Range("K55").Value = 1307.25
$ws.Range("M55").Value = -1134.25

$ws.Range("H68").Value = 7066.524
$ws.Range("I68").Value = 5249.75
$ws.Range("J68").Value = 8184.5386
$ws.Range("K68").Value = 5249.75
$ws.Range("L68").Value = 8184.5386
$ws.Range("M68").Value = -4500.75
$ws.Range("N68").Value = -9682.5386

$ws.Range("H71").Value = 7066.524
$ws.Range("I71").Value = 5249.75
$ws.Range("J71").Value = 8184.5386
$ws.Range("K71").Value = 26248.75
$ws.Range("L71").Value = 40922.693
$ws.Range("M71").Value = -22504.75
$ws.Range("N71").Value = -48410.693

$ws.Range("H132").Value = 3886.2727
$ws.Range("I132").Value = 4249.5
$ws.Range("J132").Value = 3805.5557
$ws.Range("K132").Value = 12748.5
$ws.Range("L132").Value = 11416.6671
$ws.Range("M132").Value = -10218.5
$ws.Range("N132").Value = -16476.6671

$ws.Range("H136").Value = 18333.53
$ws.Range("I136").Value = 4138.1113
$ws.Range("J136").Value = 34303.375
$ws.Range("K136").Value = 12414.3339
$ws.Range("L136").Value = 102910.125
$ws.Range("M136").Value = -9864.333899999998
$ws.Range("N136").Value = -108010.125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 53661
$ws.Range("J46").Value = 53661
$ws.Range("L46").Value = 53661
$ws.Range("N46").Value = -54123

$ws.Range("H81").Value = 2171.8125
$ws.Range("I81").Value = 1167
$ws.Range("K81").Value = 2334
$ws.Range("M81").Value = -1273

$ws.Range("H84").Value = 2171.8125
$ws.Range("I84").Value = 1167
$ws.Range("K84").Value = 11670
$ws.Range("M84").Value = -6366

$ws.Range("H100").Value = 827.8
$ws.Range("I100").Value = 457.57144
$ws.Range("K100").Value = 915.14288
$ws.Range("M100").Value = -374.14288

$ws.Range("H107").Value = 724.8182
$ws.Range("I107").Value = 640
$ws.Range("J107").Value = 951
$ws.Range("K107").Value = 1920
$ws.Range("L107").Value = 2853
$ws.Range("M107").Value = 0
$ws.Range("N107").Value = -6693

$ws.Range("H113").Value = 1036.619
$ws.Range("I113").Value = 672.0833
$ws.Range("K113").Value = 2016.2499
$ws.Range("M113").Value = 153.7501

$ws.Range("H132").Value = 5328.1016
$ws.Range("I132").Value = 5661.7
$ws.Range("J132").Value = 3474.7778
$ws.Range("K132").Value = 16985.1
$ws.Range("L132").Value = 10424.3334
$ws.Range("M132").Value = -14455.1
$ws.Range("N132").Value = -15484.3334

$ws.Range("H134").Value = 53661
$ws.Range("J134").Value = 53661
$ws.Range("L134").Value = 160983
$ws.Range("N134").Value = -166053

$ws.Range("H136").Value = 4222.643
$ws.Range("I136").Value = 2657.2593
$ws.Range("K136").Value = 7971.777900000001
$ws.Range("M136").Value = -5421.777900000001
